# feat: #108 Re-generate the statistics with the fixed minutes and seconds
# formatting in the haul fields ("Квала" column I: "H ч. M мин. S сек.").
#
# Zero-pad single-digit minutes/seconds to two digits, e.g.
#   "9 ч. 5 мин. 21 сек."  -> "9 ч. 05 мин. 21 сек."
#   "0 ч. 31 мин. 0 сек."  -> "0 ч. 31 мин. 00 сек."
# Hours are left untouched (they are not zero-padded).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
            $hours = $matches[1]
            $minutes = $matches[2]
            $seconds = $matches[3]

            if ($minutes.Length -lt 2 -or $seconds.Length -lt 2) {
                $paddedMinutes = $minutes.PadLeft(2, '0')
                $paddedSeconds = $seconds.PadLeft(2, '0')
                $newVal = "$hours ч. $paddedMinutes мин. $paddedSeconds сек."
                $cell.Value = $newVal
            }
        }
    }
}

